$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header Q1: RequireMedical -> MedicareOption
$ws.Range("Q1").Value = "MedicareOption"

# Change Q2 value from FALSE to Hide
$ws.Range("Q2").Value = "Hide"

# Update selection to Q3
$ws.Range("Q3").Select()
